$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.730.80'
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").Value = '2.660.29'
$ws.Range("E3").Value = '  -0.65%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.65'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.86%  '
$ws.Range("E7").Value = '  +3.78%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.126'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.45%  '
$ws.Range("E10").Value = '  -0.45%  '
$ws.Range("E11").Value = '  -0.27%  '
$ws.Range("E12").Value = '  +1.14%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.15'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.48%  '
$ws.Range("E14").Value = '  -1.33%  '
$ws.Range("D15").Value = '3.141.65'
$ws.Range("E15").Value = '  -0.53%  '
$ws.Range("D16").Value = '65.641.43'
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").Value = '2.674.60'
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.54'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.94%  '
$ws.Range("E19").Value = '  -0.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '354.16'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.46'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.79'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("B24").Value = 'PEPE'
$ws.Range("C24").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000113'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.54%  '
$ws.Range("B25").Value = 'SuiNetwork'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.78'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.74'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.87%  '
$ws.Range("E27").Value = '  +1.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '562.87'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.10'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.60%  '
$ws.Range("E30").Value = '  -2.24%  '
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.14'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.15%  '
$ws.Range("E33").Value = '  +2.95%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.70'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.50'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.01%  '
$ws.Range("E36").Value = '  -0.20%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.58'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.11%  '
$ws.Range("E38").Value = '  +1.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.998'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '153.96'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.33%  '
$ws.Range("E41").Value = '  +7.29%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '161.43'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.49%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.09'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.10%  '
$ws.Range("E44").Value = '  +1.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '23.40'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.43%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.645'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0258'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.80'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.29%  '
$ws.Range("E50").Value = '  -6.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.815'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.25%  '
